$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 97, shifting existing rows 97-122 down to 99-124.
$ws.Rows("97:98").Insert()

# Populate new row 97 with data (Ciruela / Vega Central Mapocho de Santiago weekly entry).
$ws.Cells.Item(97, 1).Value = 9
$ws.Cells.Item(97, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(97, 3).Value = "Metropolitana"
$ws.Cells.Item(97, 4).Value = 44637
$ws.Cells.Item(97, 5).Value = 13
$ws.Cells.Item(97, 6).Value = "Fruta"
$ws.Cells.Item(97, 7).Value = 100103
$ws.Cells.Item(97, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(97, 9).Value = 100103002
$ws.Cells.Item(97, 10).Value = "Ciruela"
$ws.Cells.Item(97, 11).Value = "Angeleno"
$ws.Cells.Item(97, 12).Value = "Especial"
$ws.Cells.Item(97, 13).Value = 270
$ws.Cells.Item(97, 14).Value = 7500
$ws.Cells.Item(97, 15).Value = 7500
$ws.Cells.Item(97, 16).Value = 7500
$ws.Cells.Item(97, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(97, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(97, 19).Value = 500
$ws.Cells.Item(97, 20).Value = 15

# Populate new row 98 with data.
$ws.Cells.Item(98, 1).Value = 9
$ws.Cells.Item(98, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(98, 3).Value = "Metropolitana"
$ws.Cells.Item(98, 4).Value = 44637
$ws.Cells.Item(98, 5).Value = 13
$ws.Cells.Item(98, 6).Value = "Fruta"
$ws.Cells.Item(98, 7).Value = 100103
$ws.Cells.Item(98, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(98, 9).Value = 100103002
$ws.Cells.Item(98, 10).Value = "Ciruela"
$ws.Cells.Item(98, 11).Value = "Angeleno"
$ws.Cells.Item(98, 12).Value = "Primera"
$ws.Cells.Item(98, 13).Value = 250
$ws.Cells.Item(98, 14).Value = 6000
$ws.Cells.Item(98, 15).Value = 6000
$ws.Cells.Item(98, 16).Value = 6000
$ws.Cells.Item(98, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(98, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(98, 19).Value = 400
$ws.Cells.Item(98, 20).Value = 15
